$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the structure of the last existing row (including the blank
# K/M placeholder cells) onto the new row, then overwrite with this
# flight's data.
$ws.Range("A9:M9").Copy($ws.Range("A10:M10"))

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Monday, Jan 09"
$ws.Range("C10").Value = "2:02 PM"
$ws.Range("D10").Value = "UNKNOWN"
$ws.Range("E10").Value = "Palma de"
$ws.Range("F10").Value = "(PMI)"
$ws.Range("G10").Value = "AMC Aviation "
$ws.Range("H10").Value = "C25B"
$ws.Range("I10").Value = "(SP-KOW)"
$ws.Range("J10").Value = "1:53 PM"
$ws.Range("L10").Value = "0 hours, -9 minutes"
